$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting rows 54:137 down to 55:138
$ws.Rows("54:54").Insert()

# Populate the newly inserted row 54 with the new data record
$ws.Cells.Item(54, 1).Value = 10
$ws.Cells.Item(54, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(54, 3).Value = "La Araucanía"
$ws.Cells.Item(54, 4).Value = 44483
$ws.Cells.Item(54, 5).Value = 9
$ws.Cells.Item(54, 6).Value = 100112013
$ws.Cells.Item(54, 7).Value = "Alcachofa"
$ws.Cells.Item(54, 8).Value = "Española"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 150
$ws.Cells.Item(54, 11).Value = 14000
$ws.Cells.Item(54, 12).Value = 14000
$ws.Cells.Item(54, 13).Value = 14000
$ws.Cells.Item(54, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(54, 15).Value = "Región Metropolitana"
$ws.Cells.Item(54, 16).Value = 467
$ws.Cells.Item(54, 17).Value = 30
$ws.Cells.Item(54, 18).Value = "Hortaliza"
